# New crime data collected - weekly CompStat update (cs-en-us-102pct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-range banner) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Cells that flip from numeric to the "0" / "***.*" text markers ---
# (set the text first - using a leading apostrophe for the numeric-looking
# "0" so it is stored as text - then paste the formatting from an existing
# style-14 text cell on top so the cell style matches the target exactly)
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"

$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "***.*"

$ws.Range("C27").Value = "'0"

$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Row 15 (Rape) remaining numeric updates ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("L15").Value = 18.75

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 87.5
$ws.Range("I16").Value = 176
$ws.Range("J16").Value = 98
$ws.Range("K16").Value = 79.591836734693
$ws.Range("L16").Value = 55.752212389380
$ws.Range("M16").Value = -26.359832635983
$ws.Range("N16").Value = -82.608695652173

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 324
$ws.Range("J17").Value = 234
$ws.Range("K17").Value = 38.461538461538
$ws.Range("L17").Value = 49.308755760368
$ws.Range("M17").Value = 66.153846153846
$ws.Range("N17").Value = 1.567398119122

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 46.153846153846
$ws.Range("L18").Value = 5.555555555555
$ws.Range("M18").Value = -55.218855218855
$ws.Range("N18").Value = -89.169381107491

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 11.538461538461
$ws.Range("I19").Value = 322
$ws.Range("J19").Value = 221
$ws.Range("K19").Value = 45.701357466063
$ws.Range("L19").Value = 33.057851239669
$ws.Range("M19").Value = 8.417508417508
$ws.Range("N19").Value = -34.419551934826

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 61.538461538461
$ws.Range("I20").Value = 220
$ws.Range("J20").Value = 196
$ws.Range("K20").Value = 12.244897959183
$ws.Range("L20").Value = 70.542635658914
$ws.Range("M20").Value = -10.931174089068
$ws.Range("N20").Value = -92.577597840755

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 13.636363636363
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = 28.089887640449
$ws.Range("I21").Value = 1199
$ws.Range("J21").Value = 865
$ws.Range("K21").Value = 38.612716763005
$ws.Range("L21").Value = 41.391509433962
$ws.Range("M21").Value = -7.126258714175
$ws.Range("N21").Value = -80.162144275314

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -10
$ws.Range("F24").Value = 130
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 21.495327102803
$ws.Range("I24").Value = 1410
$ws.Range("J24").Value = 862
$ws.Range("K24").Value = 63.573085846867
$ws.Range("L24").Value = 74.937965260545
$ws.Range("M24").Value = 105.839416058394

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 35.294117647058
$ws.Range("I25").Value = 471
$ws.Range("J25").Value = 435
$ws.Range("K25").Value = 8.275862068965
$ws.Range("L25").Value = 25.935828877005
$ws.Range("M25").Value = -13.736263736263

# --- Row 26 (UCR Rape*) remaining numeric updates ---
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("L26").Value = 14.285714285714

# --- Row 27 (Other Sex Crimes) remaining numeric updates ---
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = -26.666666666666
$ws.Range("L27").Value = 0

"Applied weekly crime-data refresh"
